$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update filename text in A2
$ws.Range("A2").Value = "Waves_011.txt"

# Update numeric statistic values in row 2
$ws.Range("D2").Value = 50
$ws.Range("E2").Value = 27
$ws.Range("F2").Value = 413.07
$ws.Range("G2").Value = 10.27
$ws.Range("H2").Value = 14.67
$ws.Range("K2").Value = 9.75
$ws.Range("L2").Value = 1.47
$ws.Range("M2").Value = 0.04
$ws.Range("N2").Value = 8.59
$ws.Range("O2").Value = 1.33
$ws.Range("P2").Value = 0.04
$ws.Range("Q2").Value = 142.93
$ws.Range("R2").Value = 29.44
$ws.Range("S2").Value = 0.8100000000000001
$ws.Range("T2").Value = 7.69
$ws.Range("W2").Value = 326.77
$ws.Range("X2").Value = 37.55
$ws.Range("Y2").Value = 1.03
$ws.Range("Z2").Value = 21.88
$ws.Range("AC2").Value = 18.36
$ws.Range("AD2").Value = 2.22
$ws.Range("AF2").Value = 22.48
$ws.Range("AG2").Value = 2.07
$ws.Range("AH2").Value = 0.06
$ws.Range("AI2").Value = 25.28
$ws.Range("AJ2").Value = 1.9
$ws.Range("AK2").Value = 0.05
